# Updated handling of dot correction to be compatible with LPA-Program 1.0.0.
#
# - Rename the "Archive ID" column header to "LED Set".
# - Convert the "Channel" column from text values ("Top"/"Bot") to the
#   numeric dot-correction channel values (1/2) expected by LPA-Program 1.0.0.
# - Refresh the saved view (zoom level / active selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename header "Archive ID" -> "LED Set"
$headerCell = $ws.Cells.Item(1, 1)
if ($headerCell.Text -eq "Archive ID") {
    $headerCell.Value = "LED Set"
}

# Convert "Channel" column (column C) from text Top/Bot to numeric 1/2
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Text
    if ($val -eq "Top") {
        $cell.Value = 1
    } elseif ($val -eq "Bot") {
        $cell.Value = 2
    }
}

# Update the saved view to match the author's state (zoom 100%, E7 selected)
$excel.ActiveWindow.Zoom = 100
$ws.Range("E7").Select()
